$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Area 1: "...Project document...(.doc or .docx) or power point or PDF only.
#          No other format..." paragraph.
#
# Before:  "or " | [bookmarkStart _GoBack] "power point" | " or PDF" |
#          " only" | [bookmarkEnd _GoBack] "." | " No other format..."
# After:   "or power point" | " or PDF" | " only." | " No other format..."
#
# The _GoBack bookmark that used to wrap "power point ... only" is removed
# entirely and the runs on either side of its old start/end tags are merged.
# ---------------------------------------------------------------------------

# Known fixed offsets (characters, not counting bookmark anchors) in the
# un-edited document:
#   "or "          [1850,1853)
#   "power point"  [1853,1864)
#   " or PDF"      [1864,1871)
#   " only"        [1871,1876)
#   "."            [1876,1877)
# followed immediately by " No other format will be acceptable."

# Locate them robustly via Find so the script is resilient to minor drift.
# Each subsequent search is scoped to start right after the previous match
# so we don't accidentally latch onto an earlier, unrelated occurrence of
# the same text elsewhere in the document (e.g. "... shows up only ...").
$rFind = $d.Content
$rFind.Find.Execute("or power point", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$orStart = $rFind.Start
$ppEnd = $rFind.End

$rFind2 = $d.Range($ppEnd, $d.Content.End)
$rFind2.Find.Execute(" or PDF", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pdfStart = $rFind2.Start
$pdfEnd = $rFind2.End

$rFind3 = $d.Range($pdfEnd, $d.Content.End)
$rFind3.Find.Execute(" only", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$onlyStart = $rFind3.Start
$onlyEnd = $rFind3.End

$dotStart = $onlyEnd
$dotEnd = $dotStart + 1

# Blockers prevent the normalize-on-edit pass from merging further than we
# want: once text inside a run changes, Word merges the touched run with any
# run-neighbours that share identical formatting, cascading until it meets a
# bookmark boundary (or the paragraph edge). We plant temporary bookmarks at
# the junctions we want to protect, do the edits, then remove them.
$d.Bookmarks.Add("ZZEDITBLOCK_L1", $d.Range($orStart, $orStart)) | Out-Null
$d.Bookmarks.Add("ZZEDITBLOCK_M1", $d.Range($pdfStart, $pdfEnd)) | Out-Null
$d.Bookmarks.Add("ZZEDITBLOCK_R1", $d.Range($dotEnd, $dotEnd)) | Out-Null

# Remove the original _GoBack bookmark (both its start and end tags).
$d.Bookmarks("_GoBack").Delete()

# Merge "or " + "power point" -> "or power point".
$rMerge1 = $d.Range($orStart, $ppEnd)
$mergedText1 = $rMerge1.Text
$rMerge1.Text = $mergedText1 + [char]1
$d.Range($ppEnd, $ppEnd + 1).Text = ""

# Merge " only" + "." -> " only.".
$rMerge2 = $d.Range($onlyStart, $dotEnd)
$mergedText2 = $rMerge2.Text
$rMerge2.Text = $mergedText2 + [char]1
$d.Range($dotEnd, $dotEnd + 1).Text = ""

# Remove the temporary blockers.
$d.Bookmarks("ZZEDITBLOCK_L1").Delete()
$d.Bookmarks("ZZEDITBLOCK_M1").Delete()
$d.Bookmarks("ZZEDITBLOCK_R1").Delete()

# ---------------------------------------------------------------------------
# Area 2: "The presentation and group evaluation will be on Saturday,
#          April 11, 2019 from 8:00 am - noon." paragraph.
#
# Before: single run "from 8:00 am - noon" with no bookmark in this
#         paragraph.
# After:  a fresh _GoBack bookmark starts right after "...will be on " and
#         ends between "from 8:00 " and "am - noon" (splitting that run in
#         two).
# ---------------------------------------------------------------------------

$rOn = $d.Content
$rOn.Find.Execute("will be on ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$onEnd = $rOn.End

$rNoon = $d.Range($onEnd, $d.Content.End)
$rNoon.Find.Execute("from 8:00 ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$from800End = $rNoon.End

$d.Bookmarks.Add("_GoBack", $d.Range($onEnd, $from800End)) | Out-Null

Write-Output "done"
